$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23: 2012-10-11 (serial 41193), 2.25h effort, "Manual continued"
$ws.Range("A23").Value = 41193
$ws.Range("B23").Value = 2.25
$ws.Range("D23").Value = "Manual continued"

# Row 24: 2012-10-12 (serial 41194), 5h effort, "Manual continued"
$ws.Range("A24").Value = 41194
$ws.Range("B24").Value = 5
$ws.Range("D24").Value = "Manual continued"

# Row 25: 2012-10-15 (serial 41197), 1h effort, fix description
$ws.Range("A25").Value = 41197
$ws.Range("B25").Value = 1
$ws.Range("D25").Value = "Fix: Bad specification of ALL events - now timer events are still an OR condition"

# Update selection to match the recorded state after the edit
$ws.Range("E25").Select() | Out-Null
